$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the decorative "BlaBlaBla" helper cells that lived beside the table
$ws.Range("F8").Clear()
$ws.Range("H10").Clear()
$ws.Range("G13").Clear()

# Remove the Go-gopher logo/picture that sat above the table
$ws.Shapes.Item("Picture 2").Delete()

# Drop the blank banner rows (1-4) that used to hold the picture; this
# shifts the whole table up so row 5 becomes row 1, etc.
$ws.Range("A1:A4").EntireRow.Delete()

# Give column B (now hidden alongside A) back its normal width before
# grouping/hiding, matching the rest of the sheet's default column width
$ws.Range("B1").EntireColumn.ColumnWidth = 9.95

# Group + hide the description/detail columns (A and B), leaving an
# outline button so the reader can re-expand them
$ws.Range("A1:B1").EntireColumn.Group()
$ws.Range("A1:B1").EntireColumn.Hidden = $true

# Update the view: zoomed to 125%, scrolled so column C is the first
# visible column, selection resting on A3
$excel.ActiveWindow.Zoom = 125
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("A3").Select()
